# "Refined metadata to be additional tab"
#
# 1. Add a new "metadata" worksheet (after "data") summarising the
#    PanelApp query that produced the "data" sheet.
# 2. Refresh the "time_taken" column (F) on the "data" sheet to the
#    timestamps of the latest query run.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# ---- 1. Create the "metadata" sheet right after "data" ----
$newSheet = $wb.Worksheets.Add($null, $dataSheet)
$newSheet.Name = "metadata"

# Reuse the exact header styling (bold, thin border, centered) from the
# "data" sheet's own header row instead of rebuilding it from scratch.
$dataSheet.Range("B1").Copy() | Out-Null
$newSheet.Range("B1:G1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$newSheet.Range("A2").PasteSpecial(-4122) | Out-Null       # xlPasteFormats

$newSheet.Cells.Item(1, 2).Value = "data_name"
$newSheet.Cells.Item(1, 3).Value = "data_id"
$newSheet.Cells.Item(1, 4).Value = "data_version"
$newSheet.Cells.Item(1, 5).Value = "data_version_created"
$newSheet.Cells.Item(1, 6).Value = "panel_query_time"
$newSheet.Cells.Item(1, 7).Value = "panel_get_request"

$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).Value = "Gastrointestinal neuromuscular disease"
$newSheet.Cells.Item(2, 3).Value = 3087

# Keep "1.15" as text (not auto-converted to the number 1.15).
$versionCell = $newSheet.Cells.Item(2, 4)
$versionCell.NumberFormat = "@"
$versionCell.Value = "1.15"
$versionCell.ClearFormats()

$newSheet.Cells.Item(2, 5).Value = "2021-08-03T22:27:59.603590Z"
$newSheet.Cells.Item(2, 6).Value = "2021-10-05 14:33:52.140087"
$newSheet.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/3087/?format=json"

# Keep "data" as the active/visible sheet.
$dataSheet.Activate()

# ---- 2. Refresh time_taken (column F) on the "data" sheet ----
$ws = $dataSheet
$ws.Cells.Item(2, 6).Value = "2021-10-05 14:33:52.143926"
$ws.Cells.Item(3, 6).Value = "2021-10-05 14:33:52.143934"
$ws.Cells.Item(4, 6).Value = "2021-10-05 14:33:52.143937"
$ws.Cells.Item(5, 6).Value = "2021-10-05 14:33:52.143940"
$ws.Cells.Item(6, 6).Value = "2021-10-05 14:33:52.143942"
$ws.Cells.Item(7, 6).Value = "2021-10-05 14:33:52.143945"
$ws.Cells.Item(8, 6).Value = "2021-10-05 14:33:52.143948"
$ws.Cells.Item(9, 6).Value = "2021-10-05 14:33:52.143950"
$ws.Cells.Item(10, 6).Value = "2021-10-05 14:33:52.143953"
$ws.Cells.Item(11, 6).Value = "2021-10-05 14:33:52.143955"
$ws.Cells.Item(12, 6).Value = "2021-10-05 14:33:52.143958"
$ws.Cells.Item(13, 6).Value = "2021-10-05 14:33:52.143960"
$ws.Cells.Item(14, 6).Value = "2021-10-05 14:33:52.143963"
$ws.Cells.Item(15, 6).Value = "2021-10-05 14:33:52.143965"
$ws.Cells.Item(16, 6).Value = "2021-10-05 14:33:52.143968"
$ws.Cells.Item(17, 6).Value = "2021-10-05 14:33:52.143970"
$ws.Cells.Item(18, 6).Value = "2021-10-05 14:33:52.143973"
$ws.Cells.Item(19, 6).Value = "2021-10-05 14:33:52.143976"
$ws.Cells.Item(20, 6).Value = "2021-10-05 14:33:52.143979"
$ws.Cells.Item(21, 6).Value = "2021-10-05 14:33:52.143981"
$ws.Cells.Item(22, 6).Value = "2021-10-05 14:33:52.143984"
$ws.Cells.Item(23, 6).Value = "2021-10-05 14:33:52.143986"
$ws.Cells.Item(24, 6).Value = "2021-10-05 14:33:52.143989"
$ws.Cells.Item(25, 6).Value = "2021-10-05 14:33:52.143991"
$ws.Cells.Item(26, 6).Value = "2021-10-05 14:33:52.143994"
$ws.Cells.Item(27, 6).Value = "2021-10-05 14:33:52.143997"
$ws.Cells.Item(28, 6).Value = "2021-10-05 14:33:52.143999"
$ws.Cells.Item(29, 6).Value = "2021-10-05 14:33:52.144002"
$ws.Cells.Item(30, 6).Value = "2021-10-05 14:33:52.144004"
$ws.Cells.Item(31, 6).Value = "2021-10-05 14:33:52.144007"
$ws.Cells.Item(32, 6).Value = "2021-10-05 14:33:52.144009"
$ws.Cells.Item(33, 6).Value = "2021-10-05 14:33:52.144012"
$ws.Cells.Item(34, 6).Value = "2021-10-05 14:33:52.144014"
$ws.Cells.Item(35, 6).Value = "2021-10-05 14:33:52.144017"
$ws.Cells.Item(36, 6).Value = "2021-10-05 14:33:52.144020"
$ws.Cells.Item(37, 6).Value = "2021-10-05 14:33:52.144022"
$ws.Cells.Item(38, 6).Value = "2021-10-05 14:33:52.144025"
$ws.Cells.Item(39, 6).Value = "2021-10-05 14:33:52.144027"
$ws.Cells.Item(40, 6).Value = "2021-10-05 14:33:52.144029"

Write-Output "Workbook updated: metadata sheet added, time_taken refreshed."
